$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.633.34'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +4.21%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''1.752.57'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +4.93%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''0.9975'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  -0.11%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''247.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +3.54%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''0.9978'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -0.19%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.4818'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  +0.05%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.2710'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +3.02%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.06256'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +1.11%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''1.743.54'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +4.37%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.07118'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +1.48%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''15.98'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +7.10%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''0.6238'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +5.53%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''4.516'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +2.94%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''77.36'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +2.78%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''0.9978'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -0.19%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''26.644.11'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +4.28%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''0.9987'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -0.08%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''0.000006906'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +1.94%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''11.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +2.79%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''1.965.65'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +4.44%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''4.645'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +4.39%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''8.878'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +1.59%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''5.355'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +1.26%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  -0.43%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''15.50'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +3.00%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''1.840'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +6.52%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''1.415'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +1.56%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''107.99'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +2.96%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''4.028'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +1.11%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''3.774'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +3.07%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''0.07910'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +1.19%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '''  +8.20%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''2.609'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -0.33%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''1.004'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +5.46%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''0.6372'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +4.50%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''0.9494'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +10.59%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''114.69'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +19.48%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''2.511'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -3.15%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''1.997'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +7.35%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''1.003'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +0.32%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = '''VeChain'
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = '''0.01513'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +2.12%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = '''FraxShare'
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = '''5.738'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +18.39%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''0.3925'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +3.94%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''6.783'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +9.13%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.1207'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +7.72%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''0.05342'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +1.67%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''8.028'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +9.19%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''31.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +3.79%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.3468'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +3.77%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''51.87'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +3.76%  '
$ws.Range("E51").Style = "Normal"
